$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value into a cell while preserving it as
# TEXT (the workbook stores the "Price" column as inline strings, not
# numbers), matching the source data's inline-string typing.
function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
}

# --- Column D ("Price") numeric-text updates -------------------------------
Set-TextValue "D2"  "245.48"
Set-TextValue "D4"  "5.333"
Set-TextValue "D5"  "0.05837"
Set-TextValue "D7"  "3.362"
Set-TextValue "D8"  "0.8124"
Set-TextValue "D9"  "0.9210"
Set-TextValue "D11" "0.07360"
Set-TextValue "D12" "0.03085"
Set-TextValue "D13" "0.03074"
Set-TextValue "D14" "0.09359"
Set-TextValue "D15" "3.874"
Set-TextValue "D16" "0.001575"
Set-TextValue "D17" "0.04680"
Set-TextValue "D18" "0.0005994"
Set-TextValue "D19" "0.005959"
Set-TextValue "D20" "0.001247"
Set-TextValue "D21" "0.004689"
Set-TextValue "D22" "0.00008821"
Set-TextValue "D25" "0.3230"
Set-TextValue "D40" "0.03845"

# Row 18 ("One"/ONE) also gains a "Worstin24h" marker on its E cell
$ws.Range("E18").Value = "17OneONEWorstin24h"

# --- Rows 41-43: the three coins rotate position ---------------------------
# Row 41 used to be BKEXToken -> becomes KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006365"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 used to be CEJI -> becomes BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 used to be KickToken -> becomes CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003208"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining column D numeric-text updates --------------------------------
Set-TextValue "D44" "0.008495"
Set-TextValue "D45" "0.00005260"
Set-TextValue "D47" "0.6535"
Set-TextValue "D48" "0.001861"
